$wb = $excel.ActiveWorkbook
$odi = $wb.ActiveSheet

# --- 1. Update MATCH_CARD_LINK column on "ODI Batting" to MATCH_CODE with bare match codes ---
$odi.Range("D1").Value = "MATCH_CODE"

# Leading "'" forces the numeric-looking match code to stay a text value,
# matching the other text-typed columns in this sheet (e.g. MATCH_NUMBER).
$odi.Range("D2").Value  = "'4111"
$odi.Range("D3").Value  = "'4113"
$odi.Range("D4").Value  = "'4118"
$odi.Range("D5").Value  = "'4120"
$odi.Range("D6").Value  = "'4161"
$odi.Range("D7").Value  = "'4381"
$odi.Range("D8").Value  = "'4384"
$odi.Range("D9").Value  = "'4386"
$odi.Range("D10").Value = "'4625"
$odi.Range("D11").Value = "'4631"
$odi.Range("D12").Value = "'4632"
$odi.Range("D13").Value = "'4635"

# --- 2. Insert a new "Player Info" worksheet before "ODI Batting" ---
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$info.Range("A2").Value = "'4723"
$info.Range("B2").Value = "Michael Alexander Jones"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Off Break"
